$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Sector'
$ws.Range('B1').Value = 'Average Correlation'
$ws.Range('A2').Value = 'Multi-Utilities(18)'
$ws.Range('B2').Value = 0.6806635493997961
$ws.Range('A3').Value = 'Household Products(10)'
$ws.Range('B3').Value = 0.5551239585304982
$ws.Range('A4').Value = 'Electric Utilities(28)'
$ws.Range('B4').Value = 0.5028578579683908
$ws.Range('A5').Value = 'Gas Utilities(12)'
$ws.Range('B5').Value = 0.4996720732227125
$ws.Range('A6').Value = 'Containers & Packaging(12)'
$ws.Range('B6').Value = 0.4554211694031084
$ws.Range('A7').Value = 'Insurance(75)'
$ws.Range('B7').Value = 0.4434552474521457
$ws.Range('A8').Value = 'Equity Real Estate Investment Trusts ...(98)'
$ws.Range('B8').Value = 0.4395536947098878
$ws.Range('A9').Value = 'Road & Rail(22)'
$ws.Range('B9').Value = 0.4197651903675268
$ws.Range('A10').Value = 'Marine(15)'
$ws.Range('B10').Value = 0.3924305106678252
$ws.Range('A11').Value = 'Mortgage Real Estate Investment Trust...(16)'
$ws.Range('B11').Value = 0.3569377204730874
$ws.Range('A12').Value = 'Professional Services(35)'
$ws.Range('B12').Value = 0.3515205790689357
$ws.Range('A13').Value = 'Construction & Engineering(21)'
$ws.Range('B13').Value = 0.3438736672397285
$ws.Range('A14').Value = 'Water Utilities(13)'
$ws.Range('B14').Value = 0.3365577813442983
$ws.Range('A15').Value = 'Machinery(86)'
$ws.Range('B15').Value = 0.3273196891774619
$ws.Range('A16').Value = 'Auto Components(21)'
$ws.Range('B16').Value = 0.3195755395178446
$ws.Range('A17').Value = 'Specialty Retail(59)'
$ws.Range('B17').Value = 0.3192409114387076
$ws.Range('A18').Value = 'ETF(303)'
$ws.Range('B18').Value = 0.3133605361086169
$ws.Range('A19').Value = 'Airlines(14)'
$ws.Range('B19').Value = 0.3116040117463391
$ws.Range('A20').Value = 'Chemicals(52)'
$ws.Range('B20').Value = 0.2801910920682691
$ws.Range('A21').Value = 'Building Products(24)'
$ws.Range('B21').Value = 0.2796826390963346
$ws.Range('A22').Value = 'Consumer Finance(15)'
$ws.Range('B22').Value = 0.275584626156352
$ws.Range('A23').Value = 'Diversified Telecommunication Services(20)'
$ws.Range('B23').Value = 0.2704492994662344
$ws.Range('A24').Value = 'Diversified Consumer Services(17)'
$ws.Range('B24').Value = 0.2576480116098388
$ws.Range('A25').Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Range('B25').Value = 0.2460280766957026
$ws.Range('A26').Value = 'Aerospace & Defense(37)'
$ws.Range('B26').Value = 0.2409153622487695
$ws.Range('A27').Value = 'Banks(251)'
$ws.Range('B27').Value = 0.2380411178701261
$ws.Range('A28').Value = 'Hotels, Restaurants & Leisure(51)'
$ws.Range('B28').Value = 0.2238821841357235
$ws.Range('A29').Value = 'Media(42)'
$ws.Range('B29').Value = 0.223127485609525
$ws.Range('A30').Value = 'Commercial Services & Supplies(52)'
$ws.Range('B30').Value = 0.2230322642084072
$ws.Range('A31').Value = 'Electrical Equipment(28)'
$ws.Range('B31').Value = 0.2202056081988574
$ws.Range('A32').Value = 'Food Products(46)'
$ws.Range('B32').Value = 0.216020694915553
$ws.Range('A33').Value = 'IT Services(52)'
$ws.Range('B33').Value = 0.2133866368881704
$ws.Range('A34').Value = 'Trading Companies & Distributors(25)'
$ws.Range('B34').Value = 0.2072231145600254
$ws.Range('A35').Value = 'Capital Markets(76)'
$ws.Range('B35').Value = 0.1941875669495479
$ws.Range('A36').Value = 'Household Durables(39)'
$ws.Range('B36').Value = 0.1881941681680455
$ws.Range('A37').Value = 'Metals & Mining(106)'
$ws.Range('B37').Value = 0.182221272164292
$ws.Range('A38').Value = 'Health Care Providers & Services(47)'
$ws.Range('B38').Value = 0.1715759940980815
$ws.Range('A39').Value = 'Semiconductors & Semiconductor Equipment(70)'
$ws.Range('B39').Value = 0.1697578959724912
$ws.Range('A40').Value = 'Health Care Equipment & Supplies(86)'
$ws.Range('B40').Value = 0.1581208036341931
$ws.Range('A41').Value = 'Electronic Equipment, Instruments & C...(78)'
$ws.Range('B41').Value = 0.1429583227492184
$ws.Range('A42').Value = 'Communications Equipment(45)'
$ws.Range('B42').Value = 0.1408272132494436
$ws.Range('A43').Value = 'Energy Equipment & Services(38)'
$ws.Range('B43').Value = 0.1390750693558992
$ws.Range('A44').Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range('B44').Value = 0.1299470991891688
$ws.Range('A45').Value = 'Software(70)'
$ws.Range('B45').Value = 0.1292072943183608
$ws.Range('A46').Value = 'Pharmaceuticals(53)'
$ws.Range('B46').Value = 0.1272227864674263
$ws.Range('A47').Value = 'Biotechnology(128)'
$ws.Range('B47').Value = 0.1240848670924985
$ws.Range('A48').Value = 'Oil, Gas & Consumable Fuels(125)'
$ws.Range('B48').Value = 0.0934684943417734
